$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.88%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.679"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08069"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.65%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.030"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.96%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.745"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.45%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.37%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.97%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9212"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.51%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1258"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.40%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1944"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.25%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.312"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-7.91%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09398"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.63%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03707"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "7.67%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1055"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001297"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.63%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006242"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.33%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.363"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.13%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3476"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.44%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1417"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.95%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2656"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.12%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04434"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.07%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.19%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.29%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001243"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "13.88%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02863"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "15.56%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05464"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.92%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007789"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009989"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "12.56%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1418"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002233"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.19%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01188"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "13.21%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006782"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.82%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002285"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003020"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-13.73%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
